$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestsAndResult")

# --- Update measured "Output Length" values in column K (rows 7-30). ---
# The compression ratio in column L (L = K/J) is a formula and will
# recalculate automatically once the underlying K value changes.
$newOutputLengths = @{
    7  = 2060668
    8  = 1835418
    9  = 1613235
    10 = 1408475
    11 = 1117127
    12 = 897988
    13 = 702773
    14 = 505931
    15 = 363405
    16 = 224246
    17 = 117883
    18 = 23230
    19 = 693501
    20 = 622227
    21 = 520138
    22 = 451583
    23 = 376354
    24 = 302900
    25 = 228207
    26 = 179230
    27 = 126439
    28 = 83295
    29 = 43185
    30 = 9218
}

foreach ($row in $newOutputLengths.Keys) {
    $ws.Cells.Item($row, 11).Value = $newOutputLengths[$row]
}

# --- New (still empty) columns L/M reserved for the extra compression ---
# --- stages mentioned in the commit message (codewords / GZip). Only  ---
# --- M7 actually received formatting in this pass: give it the same  ---
# --- percentage number style already used throughout column L.       ---
$ws.Range("L8").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Reasonable explicit widths for the new L/M columns (close to the
# author's manually-resized widths).
$ws.Columns.Item(12).ColumnWidth = 10.5
$ws.Columns.Item(13).ColumnWidth = 10.6666666666667

# --- Cursor / selection moved from P4 to M6 while reviewing the new columns ---
$ws.Range("M6").Select()
